$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.355932831764221
$ws.Range("B1").Value = 1.508608937263489
$ws.Range("C1").Value = 1.824611067771912
$ws.Range("D1").Value = 2.808193922042847
$ws.Range("E1").Value = 15
